$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# root_1 (row 2)
$ws.Range("B2").Value = -0.5605056207
$ws.Range("C2").Value = -224.7368393
$ws.Range("D2").Value = -225.29734492
$ws.Range("E2").Value = -224.5121314331

# root_2 (row 3)
$ws.Range("B3").Value = -0.5691547179999999
$ws.Range("C3").Value = -224.66295603
$ws.Range("D3").Value = -225.23211075
$ws.Range("E3").Value = -224.5121314331

# root_3 (row 4)
$ws.Range("B4").Value = -0.5722913108
$ws.Range("C4").Value = -224.6457806
$ws.Range("D4").Value = -225.21807191
$ws.Range("E4").Value = -224.5121314331
